$d = $word.ActiveDocument

# Locate the run of text that needs to be split up and have the new
# "(src/des_extended.v)" reference inserted into the middle of it.
$target = $d.Content
$found = $target.Find.Execute(
    "changes the trigger condition for the Trojan. Once the board is more than 45-degrees tilted, the trojan will be triggered. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -eq $true) {
    $matchStart = $target.Start

    # Text segments, matching the run boundaries in the target document:
    #   "changes the trigger condition for the Trojan"
    #   " (src/"
    #   "des_extended.v"
    #   ")"
    #   ". Once the board is more than 45-degrees tilted, the trojan will be triggered. "
    $seg1 = "changes the trigger condition for the Trojan"
    $seg2 = " (src/"
    $seg3 = "des_extended.v"
    $seg4 = ")"

    # Boundary right after "...for the Trojan" (before the new text).
    $b1 = $matchStart + $seg1.Length

    # Insert the whole new parenthetical as a single block right after the
    # boundary; it will initially share a run with its neighbours.
    $insertRange = $d.Range($b1, $b1)
    $insertRange.InsertAfter($seg2 + $seg3 + $seg4)

    $b2 = $b1 + $seg2.Length
    $b3 = $b2 + $seg3.Length
    $b4 = $b3 + $seg4.Length

    # Dropping a bookmark exactly on a boundary forces the run(s) touching
    # that point to split there; deleting the bookmark afterwards removes
    # the bookmark markup but leaves the run split in place. Do this for
    # every boundary we need - including the original boundary at
    # $matchStart, which would otherwise get merged into the new text -
    # so the final run layout matches "Attached is the Verilog code that "
    # | "changes the trigger condition for the Trojan" | " (src/" |
    # "des_extended.v" | ")" | ". Once the board is more than ..." .
    $d.Bookmarks.Add("zzSplit0", $d.Range($matchStart, $matchStart)) | Out-Null
    $d.Bookmarks.Add("zzSplit1", $d.Range($b1, $b1)) | Out-Null
    $d.Bookmarks.Add("zzSplit2", $d.Range($b2, $b2)) | Out-Null
    $d.Bookmarks.Add("zzSplit3", $d.Range($b3, $b3)) | Out-Null
    $d.Bookmarks.Add("zzSplit4", $d.Range($b4, $b4)) | Out-Null

    $d.Bookmarks("zzSplit0").Delete()
    $d.Bookmarks("zzSplit1").Delete()
    $d.Bookmarks("zzSplit2").Delete()
    $d.Bookmarks("zzSplit3").Delete()
    $d.Bookmarks("zzSplit4").Delete()

    Write-Host "Inserted source reference into the Trojan-trigger paragraph."
} else {
    Write-Host "WARNING: target sentence not found; no changes made."
}
